$wb = $excel.ActiveWorkbook

# --- 总计 (summary) sheet: add a new "2022-Q3" row, push old "2022-Q2" row down ---
$wsTotal = $wb.Worksheets.Item(1)

# Copy row 2 (currently the 2022-Q2 summary row) down to row 3, keeping formats
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Range("A3").Value = 1

# Turn row 2 into the new 2022-Q3 summary row
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("D2").Value = 1.11

# --- Fund-holding detail sheets ---
$wsQ2Old = $wb.Worksheets.Item(2)

# Duplicate the existing "2022-Q2" sheet right after itself; the duplicate keeps the
# original Q2 fund-holding data untouched and becomes the new "2022-Q2" sheet.
$wsQ2Old.Copy($null, $wsQ2Old)
$wsQ2New = $wb.Worksheets.Item(3)

# The original sheet object now becomes "2022-Q3" and gets overwritten with the
# new quarter's fund-holding data. Rename it first so the duplicate can reclaim
# the "2022-Q2" name.
$wsQ3 = $wsQ2Old
$wsQ3.Name = "2022-Q3"
$wsQ2New.Name = "2022-Q2"

# Row 2 — fund 007497
$wsQ3.Range("B2:G2").NumberFormat = "@"
$wsQ3.Range("B2").Value = "007497"
$wsQ3.Range("C2").Value = "中庚价值灵动灵活配置混合"
$wsQ3.Range("D2").Value = "36.46"
$wsQ3.Range("E2").Value = "89.30"
$wsQ3.Range("F2").Value = "3.02"
$wsQ3.Range("G2").Value = "1.1011"
$wsQ3.Range("H2").Value = 6
$wsQ3.Range("B2:G2").Style = "Normal"

# Row 3 — fund 010404
$wsQ3.Range("B3:G3").NumberFormat = "@"
$wsQ3.Range("B3").Value = "010404"
$wsQ3.Range("C3").Value = "博道盛利6个月持有期混合"
$wsQ3.Range("D3").Value = "1.10"
$wsQ3.Range("E3").Value = "41.15"
$wsQ3.Range("F3").Value = "0.36"
$wsQ3.Range("G3").Value = "0.0040"
$wsQ3.Range("H3").Value = 10
$wsQ3.Range("B3:G3").Style = "Normal"
